$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.449.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +1.39%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.675.90'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +2.11%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = '''219.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +2.37%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''0.5318'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +1.72%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.04%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.2703'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +3.83%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.06398'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +1.45%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''21.83'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +5.66%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.07801'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +1.70%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''1.684.97'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +2.70%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''4.510'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +2.26%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''0.5581'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +0.83%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''0.0₅8348'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +1.93%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''65.64'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +0.97%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''26.485.42'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +1.53%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = '''  -0.02%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''4.790'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +1.86%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''193.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +2.29%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = '''  +0.95%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''6.318'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +2.39%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  +0.10%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''0.1279'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +6.04%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''140.28'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -3.44%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''7.405'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -0.01%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  +2.88%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  +4.17%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''0.06276'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +5.71%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''1.284'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +2.21%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''3.611'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +5.06%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''3.457'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +1.69%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''1.691'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +2.51%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''1.011'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +2.91%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.6150'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +9.18%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  +1.33%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''2.781'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +0.96%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = '''  +0.85%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''6.134'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +7.76%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''1.097.94'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +7.36%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.8629'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +1.31%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D43").Value = '''100.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +0.39%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''1.820.45'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +1.73%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.0₈112'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +5.36%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''58.66'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +5.14%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''8.162'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +1.42%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''1.000'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -0.10%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''0.05199'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +0.99%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''1.475'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +6.67%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  +2.07%  '
$ws.Range("E51").Style = "Normal"
